$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swap (rows 41-42) ---
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

# --- Price column (D) updates; keep as text like the source data ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.707.81'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.600.75'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.49'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.513'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.65'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.825.94'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.584.70'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.04'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.01'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₃0737'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '210.09'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.15'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.98'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.08'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.34'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0509'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.15'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.26'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.97'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.287.50'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.49'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.829'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.785'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.19'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.78'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.738.04'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.47'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0516'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.42'

# --- Volume(1h) column (E) updates ---
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("E20").Value = '  +2.21%  '
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("E22").Value = '  -3.21%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("E31").Value = '  +0.47%  '
$ws.Range("E32").Value = '  +0.42%  '
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -2.95%  '
$ws.Range("E37").Value = '  +10.67%  '
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("E40").Value = '  -1.91%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  -0.54%  '
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("E46").Value = '  -1.43%  '
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("E51").Value = '  +0.90%  '
